# Negative-case validation data for the Facebook login credentials test sheet.
# Replaces the email and first name used in the single data row, keeping the
# other two fields (password, last name) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update firstName (column C, row 2) before the email so that the new
# shared-string entries are appended in the same order as the target file:
# "Lucy" first, then the new email address.
$ws.Cells.Item(2, 3).Value2 = "Lucy"

# Update email (column A, row 2) to the new negative-test address.
$ws.Cells.Item(2, 1).Value2 = "testing.automate.facebook@gmail.com"

# Move the active selection to F2, matching the saved view state.
[void]$ws.Range("F2").Select()
